# Apply the BD.xlsx edit: rename sheet, add documentation rows for the
# "users" and "tiendaDB" (items) MongoDB collections, bold headers, resize
# column C, and set page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet from "Hoja1" to "BD"
$ws.Name = "BD"

# 2) Insert the "BD: tiendaDB" / "Collección: items" documentation block
#    first (2 new rows above the existing header), matching the order the
#    shared strings were originally authored in.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "BD: tiendaDB"
$ws.Range("A1").Font.Bold = $true
$ws.Range("C1").Value = "Collección: items"
$ws.Range("C1").Font.Bold = $true

# 3) Insert the "BD: users" / "Collección: users" documentation block above
#    that (3 new rows), pushing the tiendaDB block down to row 4.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "BD: users"
$ws.Range("A1").Font.Bold = $true
$ws.Range("C1").Value = "Collección: users"
$ws.Range("C1").Font.Bold = $true
$ws.Range("A2").Value = 'db.items.insert({mail: "next.user1@nextu.com.co", user: "usuario uno", fecNacimiento: "1980-01-01", pass: "123"})'

# 4) Widen column C to fit the new "Collección: ..." labels.
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

# 5) Page setup: letter/A4-style paper (9) and portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
